$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PC_YeuCauLan1")
$ws.Activate()

# Update the assignment cells for YCL 1 & 2 (rows 16 & 17, column F)
$ws.Range("F16").Value = "update"
$ws.Range("F17").Value = "update"

# Move the viewport / selection to match the author's saved view
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 2
$ws.Range("F17").Select()
